# Generate Report for Handback
# Marks the localization entries as handed back (in sync with en-US),
# stamps the handback datetime, and records the latest target/handback
# file links for each locale sheet (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---- Overview sheet: Status columns (zh-cn / de-de) for both rows ----
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("B2").Value = $statusHandedBack
$ovw.Range("C2").Value = $statusHandedBack
$ovw.Range("B3").Value = $statusHandedBack
$ovw.Range("C3").Value = $statusHandedBack

function Update-LocaleSheet($SheetName, $HandbackDateTime, $Row2MdUrl, $Row2XlfUrl, $Row3MdUrl, $Row3XlfUrl) {
    $ws = $wb.Worksheets.Item($SheetName)

    # Status -> Handed back
    $ws.Range("C2").Value = $statusHandedBack
    $ws.Range("C3").Value = $statusHandedBack

    # Latest Handback DateTime (was the 0001-01-01 placeholder)
    $ws.Range("H2").Value = $HandbackDateTime
    $ws.Range("H3").Value = $HandbackDateTime

    $mdName2 = $ws.Range("A2").Text
    $xlfName2 = $ws.Range("D2").Text
    $mdName3 = $ws.Range("A3").Text
    $xlfName3 = $ws.Range("D3").Text

    # Latest Target File (F) / Latest Handback File (G) -- row 2
    $ws.Hyperlinks.Add($ws.Range("F2"), $Row2MdUrl, "", "", $mdName2)
    $ws.Hyperlinks.Add($ws.Range("G2"), $Row2XlfUrl, "", "", $xlfName2)

    # Latest Target File (F) / Latest Handback File (G) -- row 3
    $ws.Hyperlinks.Add($ws.Range("F3"), $Row3MdUrl, "", "", $mdName3)
    $ws.Hyperlinks.Add($ws.Range("G3"), $Row3XlfUrl, "", "", $xlfName3)

    # Give the new link cells the same visual style as the other hyperlink cells
    $ws.Range("F2").Style = "HyperLink"
    $ws.Range("G2").Style = "HyperLink"
    $ws.Range("F3").Style = "HyperLink"
    $ws.Range("G3").Style = "HyperLink"
}

Update-LocaleSheet "zh-cn" "2016-03-18 20:40:10" `
    "https://github.com/OpenLocalizationTest/oltest/blob/0dff96cbcb6923371dd2e22541ad41d1209fafc5/e2e/7255c0d1-5d81-46fb-9ebb-9afcd157fc5d.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c497629809749b48267d488e2d1000de92fb8cf3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/7255c0d1-5d81-46fb-9ebb-9afcd157fc5d.c3e777c7be1607b4b96ce9ccda47299f3b60db51.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTest/oltest/blob/0dff96cbcb6923371dd2e22541ad41d1209fafc5/e2e/ed864e02-fef7-412f-a184-3b196d35e764.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c497629809749b48267d488e2d1000de92fb8cf3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ed864e02-fef7-412f-a184-3b196d35e764.dfc757edacb3d6deaf0cb3bc4cce670b85e27f41.zh-cn.xlf"

Update-LocaleSheet "de-de" "2016-03-18 20:40:18" `
    "https://github.com/OpenLocalizationTest/oltest/blob/0dff96cbcb6923371dd2e22541ad41d1209fafc5/e2e/7255c0d1-5d81-46fb-9ebb-9afcd157fc5d.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/369140d4000a919601ebfcc322c7b297b6210372/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/7255c0d1-5d81-46fb-9ebb-9afcd157fc5d.c3e777c7be1607b4b96ce9ccda47299f3b60db51.de-de.xlf" `
    "https://github.com/OpenLocalizationTest/oltest/blob/0dff96cbcb6923371dd2e22541ad41d1209fafc5/e2e/ed864e02-fef7-412f-a184-3b196d35e764.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/369140d4000a919601ebfcc322c7b297b6210372/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ed864e02-fef7-412f-a184-3b196d35e764.dfc757edacb3d6deaf0cb3bc4cce670b85e27f41.de-de.xlf"
